$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 612 (pushes existing rows 612..653 down to 613..654)
$ws.Rows(612).Insert()

# Fill the new row with the daily-push data point for 2026/01/11 08:00
$ws.Range("A612").NumberFormat = "@"
$ws.Range("A612").Value = "2026/01/11"
$ws.Range("A612").NumberFormat = "General"

$ws.Range("B612").Value = "日"
$ws.Range("C612").Value = 8
$ws.Range("D612").Value = 201
